# Add two new verb rows (負ける/make-lose and 勝つ/win conjugations) to the
# bottom of the verb conjugation table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (font, row height, etc.) of the last existing data
# row (113) down onto the two new rows so the new rows look consistent
# with the rest of the table.
$ws.Range("A113:H113").Copy()
$ws.Range("A114:H115").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match row height used throughout the rest of the conjugation table.
$ws.Rows.Item(114).RowHeight = $ws.Rows.Item(113).RowHeight
$ws.Rows.Item(115).RowHeight = $ws.Rows.Item(113).RowHeight

# Row 114: 負ける (makeru / to lose)
$ws.Range("A114").Value = "負ける"
$ws.Range("B114").Value = "負けて"
$ws.Range("C114").Value = "負けた"
$ws.Range("D114").Value = "負けない"
$ws.Range("E114").Value = "負けます"
$ws.Range("F114").Value = "負けよう"
$ws.Range("G114").Value = "負けられる"
$ws.Range("H114").Value = "負けられる"

# Row 115: 勝つ (katsu / to win)
$ws.Range("A115").Value = "勝つ"
$ws.Range("B115").Value = "勝って"
$ws.Range("C115").Value = "勝った"
$ws.Range("D115").Value = "勝たない"
$ws.Range("E115").Value = "勝ちます"
$ws.Range("F115").Value = "勝とう"
$ws.Range("G115").Value = "勝てる"
$ws.Range("H115").Value = "勝たれる"

# Match the saved cursor/selection position from the authored edit: the
# user ended editing at the newly added last row.
[void]$ws.Range("A115").Select()
